$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(65, 9086, 99, 1558, 95585, 7, 7537, 7775, 965, 3668, 91)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
